# Refatorando o consolidador para modelo ETL
# Atualiza os dados de absenteismo (linhas 2-11) com os novos valores
# gerados pelo pipeline ETL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Colaborador_id, B=Colaborador_nome, C=Departamento,
#          D=Motivo_da_ausencia, E=Horas_de_ausencia, F=Data_da_ausencia, G=Salario

$rows = @(
    @{ Row=2;  A=69901; B="Dr. Caio da Mata";     C="Jurídico";          D="Problemas pessoais"; E=8; F=45088; G=7554.26 },
    @{ Row=3;  A=90430; B="Bruno Aragão";         C="Jurídico";          D="Problemas pessoais"; E=2; F=45085; G=7042.08 },
    @{ Row=4;  A=51443; B="Juliana Pereira";      C="TI";                D="Problemas pessoais"; E=3; F=45093; G=5549.76 },
    @{ Row=5;  A=8742;  B="Henrique Vieira";      C="Vendas";            D="Outros";             E=1; F=45092; G=11563.03 },
    @{ Row=6;  A=44686; B="Brenda da Luz";        C="Marketing";         D="Doença";             E=7; F=45085; G=12335.76 },
    @{ Row=7;  A=33584; B="Sr. Ryan Cunha";       C="Recursos Humanos";  D="Problemas pessoais"; E=3; F=45092; G=12226.97 },
    @{ Row=8;  A=77499; B="Sr. Nicolas Azevedo";  C="Operações";         D="Problemas pessoais"; E=6; F=45083; G=10138.74 },
    @{ Row=9;  A=72543; B="Daniel Barros";        C="Operações";         D="Problemas pessoais"; E=2; F=45103; G=8421.889999999999 },
    @{ Row=10; A=59293; B="Lucas Viana";          C="Vendas";            D="Problemas pessoais"; E=5; F=45085; G=11393.71 },
    @{ Row=11; A=94136; B="Luiz Otávio Melo";     C="Engenharia";        D="Doença";             E=2; F=45078; G=10827.15 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
